$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Regular_Section_A")
$ws.Range("D20").Value = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E20").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D21").Value = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E21").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D22").Value = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E22").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D23").Value = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"
$ws.Range("D24").Value = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"
$ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"
$ws.Range("D25").Value = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"
$ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"
$ws.Range("D26").Value = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"
$ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"
$ws.Range("D27").Value = "Mon 10:30-12:00 [C101], Wed 10:30-12:00 [C101]"
$ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"
$ws.Range("D28").Value = "Mon 10:30-12:00 [C102], Wed 10:30-12:00 [C102]"
$ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"
$ws.Range("D29").Value = "Mon 10:30-12:00 [C104], Wed 10:30-12:00 [C104]"
$ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"
$ws.Range("D30").Value = "Mon 10:30-12:00 [C202], Wed 10:30-12:00 [C202]"
$ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"
$ws.Range("D31").Value = "Mon 10:30-12:00 [C203], Wed 10:30-12:00 [C203]"
$ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C101], Thu 15:30-17:00 [C101]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C102], Thu 15:30-17:00 [C102]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C104], Thu 15:30-17:00 [C104]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [C202], Thu 15:30-17:00 [C202]"
$ws.Range("E35").Value = "Fri 14:30-15:30 [C202]"

$ws = $wb.Worksheets.Item("Regular_Section_B")
$ws.Range("D20").Value = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E20").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D21").Value = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E21").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D22").Value = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E22").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D23").Value = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"
$ws.Range("D24").Value = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"
$ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"
$ws.Range("D25").Value = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"
$ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"
$ws.Range("D26").Value = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"
$ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"
$ws.Range("D27").Value = "Mon 10:30-12:00 [C101], Wed 10:30-12:00 [C101]"
$ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"
$ws.Range("D28").Value = "Mon 10:30-12:00 [C102], Wed 10:30-12:00 [C102]"
$ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"
$ws.Range("D29").Value = "Mon 10:30-12:00 [C104], Wed 10:30-12:00 [C104]"
$ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"
$ws.Range("D30").Value = "Mon 10:30-12:00 [C202], Wed 10:30-12:00 [C202]"
$ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"
$ws.Range("D31").Value = "Mon 10:30-12:00 [C203], Wed 10:30-12:00 [C203]"
$ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C101], Thu 15:30-17:00 [C101]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C102], Thu 15:30-17:00 [C102]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C104], Thu 15:30-17:00 [C104]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [C202], Thu 15:30-17:00 [C202]"
$ws.Range("E35").Value = "Fri 14:30-15:30 [C202]"

$ws = $wb.Worksheets.Item("PreMid_Section_A")
$ws.Range("D20").Value = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E20").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D21").Value = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E21").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D22").Value = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E22").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D23").Value = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"
$ws.Range("D24").Value = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"
$ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"
$ws.Range("D25").Value = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"
$ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"
$ws.Range("D26").Value = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"
$ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"
$ws.Range("D27").Value = "Mon 10:30-12:00 [C101], Wed 10:30-12:00 [C101]"
$ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"
$ws.Range("D28").Value = "Mon 10:30-12:00 [C102], Wed 10:30-12:00 [C102]"
$ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"
$ws.Range("D29").Value = "Mon 10:30-12:00 [C104], Wed 10:30-12:00 [C104]"
$ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"
$ws.Range("D30").Value = "Mon 10:30-12:00 [C202], Wed 10:30-12:00 [C202]"
$ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"
$ws.Range("D31").Value = "Mon 10:30-12:00 [C203], Wed 10:30-12:00 [C203]"
$ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C101], Thu 15:30-17:00 [C101]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C102], Thu 15:30-17:00 [C102]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C104], Thu 15:30-17:00 [C104]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [C202], Thu 15:30-17:00 [C202]"
$ws.Range("E35").Value = "Fri 14:30-15:30 [C202]"

$ws = $wb.Worksheets.Item("PreMid_Section_B")
$ws.Range("D20").Value = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E20").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D21").Value = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E21").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D22").Value = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E22").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D23").Value = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"
$ws.Range("D24").Value = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"
$ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"
$ws.Range("D25").Value = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"
$ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"
$ws.Range("D26").Value = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"
$ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"
$ws.Range("D27").Value = "Mon 10:30-12:00 [C101], Wed 10:30-12:00 [C101]"
$ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"
$ws.Range("D28").Value = "Mon 10:30-12:00 [C102], Wed 10:30-12:00 [C102]"
$ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"
$ws.Range("D29").Value = "Mon 10:30-12:00 [C104], Wed 10:30-12:00 [C104]"
$ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"
$ws.Range("D30").Value = "Mon 10:30-12:00 [C202], Wed 10:30-12:00 [C202]"
$ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"
$ws.Range("D31").Value = "Mon 10:30-12:00 [C203], Wed 10:30-12:00 [C203]"
$ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C101], Thu 15:30-17:00 [C101]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C102], Thu 15:30-17:00 [C102]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C104], Thu 15:30-17:00 [C104]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [C202], Thu 15:30-17:00 [C202]"
$ws.Range("E35").Value = "Fri 14:30-15:30 [C202]"

$ws = $wb.Worksheets.Item("PostMid_Section_A")
$ws.Range("D20").Value = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E20").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D21").Value = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E21").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D22").Value = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E22").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D23").Value = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"
$ws.Range("D24").Value = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"
$ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"
$ws.Range("D25").Value = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"
$ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"
$ws.Range("D26").Value = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"
$ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"
$ws.Range("D27").Value = "Mon 10:30-12:00 [C101], Wed 10:30-12:00 [C101]"
$ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"
$ws.Range("D28").Value = "Mon 10:30-12:00 [C102], Wed 10:30-12:00 [C102]"
$ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"
$ws.Range("D29").Value = "Mon 10:30-12:00 [C104], Wed 10:30-12:00 [C104]"
$ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"
$ws.Range("D30").Value = "Mon 10:30-12:00 [C202], Wed 10:30-12:00 [C202]"
$ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"
$ws.Range("D31").Value = "Mon 10:30-12:00 [C203], Wed 10:30-12:00 [C203]"
$ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C101], Thu 15:30-17:00 [C101]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C102], Thu 15:30-17:00 [C102]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C104], Thu 15:30-17:00 [C104]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [C202], Thu 15:30-17:00 [C202]"
$ws.Range("E35").Value = "Fri 14:30-15:30 [C202]"

$ws = $wb.Worksheets.Item("PostMid_Section_B")
$ws.Range("D20").Value = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E20").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D21").Value = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E21").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D22").Value = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E22").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D23").Value = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"
$ws.Range("D24").Value = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"
$ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"
$ws.Range("D25").Value = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"
$ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"
$ws.Range("D26").Value = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"
$ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"
$ws.Range("D27").Value = "Mon 10:30-12:00 [C101], Wed 10:30-12:00 [C101]"
$ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"
$ws.Range("D28").Value = "Mon 10:30-12:00 [C102], Wed 10:30-12:00 [C102]"
$ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"
$ws.Range("D29").Value = "Mon 10:30-12:00 [C104], Wed 10:30-12:00 [C104]"
$ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"
$ws.Range("D30").Value = "Mon 10:30-12:00 [C202], Wed 10:30-12:00 [C202]"
$ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"
$ws.Range("D31").Value = "Mon 10:30-12:00 [C203], Wed 10:30-12:00 [C203]"
$ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C101], Thu 15:30-17:00 [C101]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C102], Thu 15:30-17:00 [C102]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C104], Thu 15:30-17:00 [C104]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [C202], Thu 15:30-17:00 [C202]"
$ws.Range("E35").Value = "Fri 14:30-15:30 [C202]"
